$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.775.00'
$ws.Cells.Item(2, 5).Value = '  +0.39%  '
$ws.Cells.Item(3, 4).Value = '1.917.12'
$ws.Cells.Item(3, 5).Value = '  +1.48%  '
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '240.90'
$ws.Cells.Item(5, 5).Value = '  -2.30%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.14%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4914'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3020'
$ws.Cells.Item(8, 5).Value = '  +2.01%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06791'
$ws.Cells.Item(9, 5).Value = '  +0.06%  '
$ws.Cells.Item(10, 4).Value = '1.923.40'
$ws.Cells.Item(10, 5).Value = '  +1.83%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '17.27'
$ws.Cells.Item(11, 5).Value = '  +0.40%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.07344'
$ws.Cells.Item(12, 5).Value = '  +1.53%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.219'
$ws.Cells.Item(13, 5).Value = '  +3.12%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '88.90'
$ws.Cells.Item(14, 5).Value = '  -2.78%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6766'
$ws.Cells.Item(15, 5).Value = '  -0.16%  '
$ws.Cells.Item(16, 4).Value = '30.756.49'
$ws.Cells.Item(16, 5).Value = '  +0.43%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.000008014'
$ws.Cells.Item(17, 5).Value = '  +0.42%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '13.65'
$ws.Cells.Item(18, 5).Value = '  +3.24%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.001'
$ws.Cells.Item(19, 5).Value = '  +0.06%  '
$ws.Cells.Item(20, 4).Value = '2.160.07'
$ws.Cells.Item(20, 5).Value = '  +1.35%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.373'
$ws.Cells.Item(21, 5).Value = '  +11.35%  '
$ws.Cells.Item(22, 2).Value = 'BinanceUSD'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.002'
$ws.Cells.Item(22, 5).Value = '  +0.10%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '201.95'
$ws.Cells.Item(23, 5).Value = '  +9.61%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.344'
$ws.Cells.Item(24, 5).Value = '  +4.72%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.728'
$ws.Cells.Item(25, 5).Value = '  +3.87%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '161.55'
$ws.Cells.Item(26, 5).Value = '  +3.73%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.85'
$ws.Cells.Item(27, 5).Value = '  -0.95%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.972'
$ws.Cells.Item(28, 5).Value = '  +3.39%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.448'
$ws.Cells.Item(29, 5).Value = '  +3.52%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.362'
$ws.Cells.Item(30, 5).Value = '  +0.85%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.09211'
$ws.Cells.Item(31, 5).Value = '  +2.08%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.123'
$ws.Cells.Item(32, 5).Value = '  +3.05%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05366'
$ws.Cells.Item(33, 5).Value = '  +3.24%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.7481'
$ws.Cells.Item(34, 5).Value = '  -0.30%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.126'
$ws.Cells.Item(35, 5).Value = '  +1.07%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.703'
$ws.Cells.Item(36, 5).Value = '  -1.45%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.01865'
$ws.Cells.Item(37, 5).Value = '  +1.62%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.730'
$ws.Cells.Item(38, 5).Value = '  +2.55%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.9306'
$ws.Cells.Item(39, 5).Value = '  -0.92%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.090'
$ws.Cells.Item(40, 5).Value = '  -2.71%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.4512'
$ws.Cells.Item(41, 5).Value = '  +1.90%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '73.44'
$ws.Cells.Item(42, 5).Value = '  +26.83%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '107.94'
$ws.Cells.Item(43, 5).Value = '  +2.16%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.963'
$ws.Cells.Item(44, 5).Value = '  +3.71%  '
$ws.Cells.Item(45, 2).Value = 'Algorand'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.1405'
$ws.Cells.Item(45, 5).Value = '  +4.86%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.003'
$ws.Cells.Item(46, 5).Value = '  +0.22%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.787'
$ws.Cells.Item(47, 5).Value = '  +2.20%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '35.97'
$ws.Cells.Item(48, 5).Value = '  +7.24%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.179'
$ws.Cells.Item(49, 5).Value = '  +5.66%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.05970'
$ws.Cells.Item(50, 5).Value = '  +2.08%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.4055'
$ws.Cells.Item(51, 5).Value = '  +3.12%  '
